$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 and IF headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style/format from the existing H1 header cell onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2 through 32: I column = 1 (constant), J column = same value as H column (IP)
for ($r = 2; $r -le 32; $r++) {
    $h = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $h
}
